# Apply updated values to the PCA_VAR and PCA_Components worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "PCA_VAR" ---
$wsVar = $wb.Worksheets.Item("PCA_VAR")

$wsVar.Range("B2").Value = 0.371
$wsVar.Range("B3").Value = 0.322
$wsVar.Range("B4").Value = 0.23
$wsVar.Range("B5").Value = 0.077

# --- Sheet "PCA_Components" ---
$wsComp = $wb.Worksheets.Item("PCA_Components")

$wsComp.Range("B2").Value = -0.271
$wsComp.Range("C2").Value = 0.858
$wsComp.Range("D2").Value = 0.249
$wsComp.Range("E2").Value = 0.357

$wsComp.Range("B3").Value = -0.635
$wsComp.Range("C3").Value = -0.431
$wsComp.Range("D3").Value = 0.631
$wsComp.Range("E3").Value = 0.114

$wsComp.Range("B4").Value = 0.116
$wsComp.Range("C4").Value = 0.259
$wsComp.Range("D4").Value = 0.448
$wsComp.Range("E4").Value = -0.848

$wsComp.Range("B5").Value = 0.714
$wsComp.Range("C5").Value = -0.1
$wsComp.Range("D5").Value = 0.583
$wsComp.Range("E5").Value = 0.375
